$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Joints")
$ws.Rows("37:37").Insert()
